$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated roster data (Player, Position, Team) for rows 2-19
$players = @(
    "Chris Paul",
    "Stephon Castle",
    "Payton Pritchard",
    "Jaylen Brown",
    "Dillon Brooks",
    "Paolo Banchero",
    "Pascal Siakam",
    "Ayo Dosunmu",
    "Nikola Jokic",
    "Isaiah Stewart",
    "Rudy Gobert",
    "Jalen Green",
    "Deni Avdija",
    "Cole Anthony",
    "Chet Holmgren",
    "Russell Westbrook",
    "Jalen Suggs",
    "Jakob Poeltl"
)

$positions = @(
    "PG",
    "PG,SG",
    "PG,SG",
    "SG,SF",
    "SG,SF,PF",
    "SF,PF",
    "SF,PF,C",
    "PG,SG,SF",
    "C",
    "PF,C",
    "C",
    "PG,SG",
    "SF,PF",
    "PG",
    "PF,C",
    "PG,SG",
    "PG,SG",
    "C"
)

$teams = @(
    "San Antonio Spurs",
    "San Antonio Spurs",
    "Boston Celtics",
    "Boston Celtics",
    "Houston Rockets",
    "Orlando Magic",
    "Indiana Pacers",
    "Chicago Bulls",
    "Denver Nuggets",
    "Detroit Pistons",
    "Minnesota Timberwolves",
    "Houston Rockets",
    "Portland Trail Blazers",
    "Orlando Magic",
    "Oklahoma City Thunder",
    "Denver Nuggets",
    "Orlando Magic",
    "Toronto Raptors"
)

for ($i = 0; $i -lt $players.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $players[$i]
    $ws.Cells.Item($row, 2).Value = $positions[$i]
    $ws.Cells.Item($row, 3).Value = $teams[$i]
}
